$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# last_edited_time (column D) bumped for the rows whose Notion pages were
# touched when reward/penalty (thưởng phạt) totals were aggregated into the
# personal report.
$ws.Range("D2").Value  = "2024-08-03T21:27:00.000Z"
$ws.Range("D3").Value  = "2024-08-03T21:27:00.000Z"
$ws.Range("D6").Value  = "2024-08-03T21:27:00.000Z"
$ws.Range("D7").Value  = "2024-08-03T21:28:00.000Z"
$ws.Range("D8").Value  = "2024-08-03T21:28:00.000Z"
$ws.Range("D11").Value = "2024-08-03T21:28:00.000Z"
$ws.Range("D13").Value = "2024-08-03T21:28:00.000Z"
